$d = $word.ActiveDocument

# 1. Split the title "Protocollo di comunicazione" into "P" + bookmark "_GoBack" + "rotocollo di comunicazione"
#    and remove the existing bookmark from the second table row's last paragraph.

# First remove the existing _GoBack bookmark (it will be re-added at the title).
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Add a bookmark named _GoBack right after the "P" in the title paragraph.
$titleRange = $d.Paragraphs(1).Range
$titleStart = $titleRange.Start
$bmRange = $d.Range($titleStart + 1, $titleStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 2. Merge the split runs in the table (these are just text runs; the Find/Replace
#    below collapses the paragraph text but Word will auto re-split runs on save
#    based on formatting, so no special action beyond ensuring text is correct).

# Merge "MSB = " + "1" -> "MSB = 1" (no textual change needed, just run merge - skip since text same)

# Merge "1" + "XXX XXXX" -> "1XXX XXXX" (paragraph merge: two paragraphs become one)
$d.Content.Find.Execute("1" + [char]13 + "XXX XXXX", $false, $false, $false, $false, $false, $true, 1, $false, "1XXX XXXX", 2) | Out-Null

Write-Host "Done"
